# Update the "Price" column (D) with refreshed quotes from the data source.
# The sheet stores these prices as plain text (not numbers), so we briefly
# mark the range as Text before writing the new values - otherwise Excel's
# automatic type detection would coerce the strings like "289.65" into
# numeric cells (and introduce floating-point noise). The style is restored
# immediately afterward so the cells keep their original (default) styling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceUpdates = [ordered]@{
    "D2"  = "289.65"
    "D3"  = "21.25"
    "D4"  = "6.470"
    "D5"  = "0.06379"
    "D6"  = "3.611"
    "D7"  = "1.578"
    "D8"  = "6.613"
    "D9"  = "0.8304"
    "D10" = "0.01429"
    "D11" = "0.1702"
    "D12" = "0.08708"
    "D13" = "0.03665"
    "D14" = "0.03216"
    "D16" = "3.706"
    "D17" = "0.001635"
    "D18" = "0.04731"
    "D19" = "0.006144"
    "D20" = "0.006296"
    "D23" = "3.770"
    "D40" = "0.04851"
    "D41" = "0.007119"
    "D43" = "0.1117"
    "D44" = "0.01150"
    "D45" = "0.00006931"
    "D47" = "0.8022"
    "D48" = "0.005551"
}

$priceRange = $ws.Range("D2:D48")
$priceRange.NumberFormat = "@"

foreach ($address in $priceUpdates.Keys) {
    $ws.Range($address).Value = $priceUpdates[$address]
}

$priceRange.Style = "Normal"

Write-Output "Updated $($priceUpdates.Count) price cells"
